$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 605.7778
$ws.Range("I18").Value = 516.6667
$ws.Range("J18").Value = 784
$ws.Range("K18").Value = 516.6667
$ws.Range("L18").Value = 784
$ws.Range("M18").Value = -232.6667
$ws.Range("N18").Value = -1352

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3587.64
$ws.Range("I3").Value = 1076.1111
$ws.Range("J3").Value = 5000.375
$ws.Range("K3").Value = 1076.1111
$ws.Range("L3").Value = 5000.375
$ws.Range("M3").Value = -961.1111000000001
$ws.Range("N3").Value = -5230.375

$ws.Range("H4").Value = 220
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = -582

$ws.Range("H13").Value = 7668
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 7668
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 7668
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -7956

$ws.Range("H32").Value = 37587.5
$ws.Range("I32").Value = 7587.8438
$ws.Range("J32").Value = 133586.4
$ws.Range("K32").Value = 7587.8438
$ws.Range("L32").Value = 133586.4
$ws.Range("M32").Value = -7300.8438
$ws.Range("N32").Value = -134160.4

$ws.Range("H80").Value = 31390.125
$ws.Range("J80").Value = 31390.125
$ws.Range("L80").Value = 31390.125
$ws.Range("N80").Value = -33386.125

$ws.Range("H82").Value = 30072.2
$ws.Range("J82").Value = 30072.2
$ws.Range("L82").Value = 30072.2
$ws.Range("N82").Value = -30794.2

$ws.Range("H83").Value = 31390.125
$ws.Range("J83").Value = 31390.125
$ws.Range("L83").Value = 94170.375
$ws.Range("N83").Value = -104154.375

$ws.Range("H85").Value = 30072.2
$ws.Range("J85").Value = 30072.2
$ws.Range("L85").Value = 30072.2
$ws.Range("N85").Value = -32568.2

$ws.Range("H132").Value = 1368.9474
$ws.Range("I132").Value = 1334.7576
$ws.Range("J132").Value = 1594.6
$ws.Range("K132").Value = 4004.2728
$ws.Range("L132").Value = 4783.799999999999
$ws.Range("M132").Value = -1474.2728
$ws.Range("N132").Value = -9843.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 6000
$ws.Range("J15").Value = 6000
$ws.Range("L15").Value = 6000
$ws.Range("N15").Value = -6454

$ws.Range("H134").Value = 1875.6333
$ws.Range("I134").Value = 1802.3793
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5407.1379
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2872.1379
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2472.678
$ws.Range("I31").Value = 2177.7273
$ws.Range("J31").Value = 2648.054
$ws.Range("K31").Value = 2177.7273
$ws.Range("L31").Value = 2648.054
$ws.Range("M31").Value = -1882.7273
$ws.Range("N31").Value = -3238.054

$ws.Range("H34").Value = 2472.678
$ws.Range("I34").Value = 2177.7273
$ws.Range("J34").Value = 2648.054
$ws.Range("K34").Value = 2177.7273
$ws.Range("L34").Value = 2648.054
$ws.Range("M34").Value = -1975.7273
$ws.Range("N34").Value = -3052.054

$ws.Range("H58").Value = 3063.9443
$ws.Range("J58").Value = 1857
$ws.Range("L58").Value = 1857
$ws.Range("N58").Value = -2263

$ws.Range("H93").Value = 16614.416
$ws.Range("I93").Value = 9721.625
$ws.Range("J93").Value = 30400
$ws.Range("K93").Value = 9721.625
$ws.Range("L93").Value = 30400
$ws.Range("M93").Value = -7849.625
$ws.Range("N93").Value = -34144

$ws.Range("H136").Value = 3063.9443
$ws.Range("J136").Value = 1857
$ws.Range("L136").Value = 5571
$ws.Range("N136").Value = -10671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1725
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 1725
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 5175
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -5851

$ws.Range("H122").Value = 2740.116
$ws.Range("I122").Value = 209
$ws.Range("J122").Value = 3504.2263
$ws.Range("K122").Value = 1881
$ws.Range("L122").Value = 31538.0367
$ws.Range("M122").Value = 569
$ws.Range("N122").Value = -36438.0367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9498.25
$ws.Range("J5").Value = 9498.25
$ws.Range("L5").Value = 9498.25
$ws.Range("N5").Value = -9722.25

$ws.Range("H103").Value = 36666.668
$ws.Range("J103").Value = 36666.668
$ws.Range("L103").Value = 36666.668
$ws.Range("N103").Value = -39010.668

$ws.Range("H122").Value = 6427.143
$ws.Range("I122").Value = 5497.5
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 16492.5
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").Value = -14042.5
$ws.Range("N122").Value = -27899.9995

$ws.Range("H126").Value = 5775.931
$ws.Range("I126").Value = 2842.1428
$ws.Range("J126").Value = 8514.134
$ws.Range("K126").Value = 8526.428400000001
$ws.Range("L126").Value = 25542.402
$ws.Range("M126").Value = -6056.428400000001
$ws.Range("N126").Value = -30482.402

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I25").Value = 5600
$ws.Range("J25").Value = 17008
$ws.Range("K25").Value = 5600
$ws.Range("L25").Value = 17008
$ws.Range("M25").Value = -5370
$ws.Range("N25").Value = -17468

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H100").Value = 3980
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3980
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3980
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -5062

$ws.Range("H132").Value = 2134.5789
$ws.Range("I132").Value = 1487.8334
$ws.Range("J132").Value = 3243.2856
$ws.Range("K132").Value = 4463.5002
$ws.Range("L132").Value = 9729.856800000001
$ws.Range("M132").Value = -1933.5002
$ws.Range("N132").Value = -14789.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4600
$ws.Range("J4").Value = 6333.3335
$ws.Range("L4").Value = 6333.3335
$ws.Range("N4").Value = -6559.3335

$ws.Range("H107").Value = 2011.1786
$ws.Range("I107").Value = 1773.1666
$ws.Range("J107").Value = 2439.6
$ws.Range("K107").Value = 5319.4998
$ws.Range("L107").Value = 7318.799999999999
$ws.Range("M107").Value = -3399.4998
$ws.Range("N107").Value = -11158.8
